$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('C摩尔-U', 'C摩尔-U', '永辉超市'),
    @('航天动力', '永辉超市', '摩尔线程'),
    @('永辉超市', '中兴通讯', '中兴通讯'),
    @('航天发展', '航天动力', '航天动力'),
    @('中兴通讯', '东百集团', '平潭发展'),
    @('永鼎股份', '永鼎股份', '航天发展'),
    @('四川金顶', '安妮股份', '顺灏股份'),
    @('东百集团', '航天发展', '安妮股份'),
    @('安妮股份', '万  科Ａ', '东百集团'),
    @('平潭发展', '雪人集团', '雪人集团'),
    @('再升科技', '平潭发展', '国晟科技'),
    @('雪人集团', '雷科防务', '四川金顶'),
    @('华菱线缆', '龙洲股份', '永鼎股份'),
    @('龙洲股份', '罗 牛 山', '银河电子'),
    @('雷科防务', '金风科技', '龙洲股份'),
    @('顺灏股份', '银河电子', '再升科技'),
    @('跃岭股份', '赛微电子', '雷科防务'),
    @('银河电子', '华菱线缆', '博纳影业'),
    @('赛微电子', '四川金顶', '赛微电子'),
    @('万  科Ａ', '顺灏股份', '茂业商业'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}